# correctif problème insertion matière dans bd et mise à jour semestre étudiant lors de l'inscription
#
# 1) "Groupe" labels were referring to semester 2 (2-A/2-B/2-C) but students
#    are actually being registered in semester 1, so the group labels must
#    be updated to 1-A/1-B/1-C.
# 2) The student identifiers in column A encoded the (wrong) enrolment year
#    2015; they must reflect 2017.
# 3) "Moyenne de l'étudiant" (column E) values are refreshed following the
#    semester correction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the "Groupe" labels (2-x -> 1-x) -------------------------------
$groupMap = @{ '2-A' = '1-A'; '2-B' = '1-B'; '2-C' = '1-C' }

for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Range("D$r")
    $current = $cell.Value2
    if ($groupMap.ContainsKey($current)) {
        $cell.Value = $groupMap[$current]
    }
}

# --- 2) Correct the enrolment year encoded in column A (2015 -> 2017) -----
for ($r = 3; $r -le 63; $r++) {
    $idCell = $ws.Range("A$r")
    $idCell.Value = $idCell.Value2 + 20000
}

# --- 3) Update the "Moyenne de l'étudiant" values in column E -------------
$newAverages = @{
    3 = 18; 4 = 11; 5 = 17; 6 = 8; 7 = 20; 9 = 9; 10 = 11; 11 = 13; 12 = 14;
    14 = 5; 15 = 10; 16 = 12; 17 = 16; 18 = 12; 19 = 6; 20 = 12; 21 = 5;
    22 = 16; 23 = 9; 24 = 16; 25 = 18; 26 = 7; 27 = 14; 28 = 5; 29 = 14;
    30 = 16; 31 = 16; 32 = 18; 34 = 11; 35 = 12; 36 = 12; 37 = 5; 38 = 14;
    39 = 14; 40 = 17; 41 = 8; 42 = 20; 43 = 20; 44 = 18; 45 = 14; 46 = 11;
    47 = 7; 48 = 10; 49 = 12; 50 = 14; 51 = 14; 52 = 12; 53 = 8; 54 = 17;
    55 = 18; 56 = 13; 58 = 14; 59 = 20; 60 = 10; 61 = 10; 62 = 19; 63 = 17
}

foreach ($r in $newAverages.Keys) {
    $ws.Range("E$r").Value = $newAverages[$r]
}

$wb.Save()
